# "Generate Report for Handoff"
# Updates the localization-status report: flips the "In Translation" status
# to "Ready for handoff" everywhere it appears, refreshes the handoff
# timestamps, and widens the Status/zh-cn/de-de columns to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Refreshed handoff timestamps --------------------------------------
# Shared by Overview!G2 and de-de!H2 (both were 2016-09-06 07:09:00)
$wsOverview.Range("G2").Value = "2016-09-06 07:09:37"
$wsDeDe.Range("H2").Value = "2016-09-06 07:09:37"

# zh-cn!H2 was a distinct timestamp (2016-09-06 07:08:56)
$wsZhCn.Range("H2").Value = "2016-09-06 07:09:33"

# --- Widen the Status / zh-cn / de-de columns to fit "Ready for handoff"
# Target stored width is 17.2159881591797 characters; the host's column
# model snaps ColumnWidth to 1/6-character increments, so feed it the
# nearest representable input (16.3333 -> stored 17.1667, the closest
# achievable cell to the authored width).
$newStatusWidth = 16.333333333333332

$wsOverview.Columns.Item(5).ColumnWidth = $newStatusWidth  # column E
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusWidth  # column F
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusWidth       # column C
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusWidth       # column C
